$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '20.570.24'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '1.477.16'
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").Value = '0.9597'
$ws.Range("E5").Value = '  +1.84%  '
$ws.Range("D6").Value = '277.43'
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("D8").Value = '0.3074'
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").Value = '1.081'
$ws.Range("E9").Value = '  +5.19%  '
$ws.Range("D10").Value = '39.36'
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").Value = '0.06648'
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '18.13'
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("D14").Value = '5.460'
$ws.Range("E14").Value = '  +2.17%  '
$ws.Range("D15").Value = '6.171'
$ws.Range("E15").Value = '  +1.75%  '
$ws.Range("D16").Value = '0.9598'
$ws.Range("E16").Value = '  +0.40%  '
$ws.Range("D17").Value = '0.00001014'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").Value = '1.475.80'
$ws.Range("E18").Value = '  +2.25%  '
$ws.Range("D19").Value = '0.05979'
$ws.Range("E19").Value = '  +5.16%  '
$ws.Range("D20").Value = '69.05'
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = '5.475'
$ws.Range("E21").Value = '  +1.61%  '
$ws.Range("E22").Value = '  +1.37%  '
$ws.Range("E23").Value = '  +2.73%  '
$ws.Range("D24").Value = '2.271'
$ws.Range("E24").Value = '  +1.55%  '
$ws.Range("D25").Value = '20.600.39'
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("D26").Value = '147.05'
$ws.Range("E26").Value = '  +4.35%  '
$ws.Range("D27").Value = '2.068'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '17.15'
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").Value = '1.637.21'
$ws.Range("E29").Value = '  +2.76%  '
$ws.Range("D30").Value = '114.75'
$ws.Range("E30").Value = '  +3.55%  '
$ws.Range("D31").Value = '3.914'
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("D32").Value = '4.926'
$ws.Range("E32").Value = '  +2.96%  '
$ws.Range("D33").Value = '0.07923'
$ws.Range("E33").Value = '  +2.87%  '
$ws.Range("D34").Value = '0.7946'
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("D35").Value = '1.188'
$ws.Range("E35").Value = '  +6.57%  '
$ws.Range("D36").Value = '1.430'
$ws.Range("E36").Value = '  -2.44%  '
$ws.Range("D37").Value = '0.05671'
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("D38").Value = '4.703'
$ws.Range("E38").Value = '  +1.25%  '
$ws.Range("D39").Value = '0.9602'
$ws.Range("E39").Value = '  +1.43%  '
$ws.Range("D40").Value = '0.02013'
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("D42").Value = '0.1839'
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("D43").Value = '7.337'
$ws.Range("E43").Value = '  +2.52%  '
$ws.Range("D44").Value = '3.513'
$ws.Range("E44").Value = '  +1.19%  '
$ws.Range("D45").Value = '0.5213'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("D46").Value = '12.04'
$ws.Range("E46").Value = '  +1.80%  '
$ws.Range("D47").Value = '119.71'
$ws.Range("E47").Value = '  +2.92%  '
$ws.Range("D48").Value = '0.5155'
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").Value = '1.808'
$ws.Range("E49").Value = '  +4.45%  '
$ws.Range("E50").Value = '  +0.42%  '
$ws.Range("D51").Value = '0.9902'
$ws.Range("E51").Value = '  +1.04%  '

# Restore default style on the price column so no stray formatting is left behind
$ws.Range("D2:D51").Style = "Normal"

